# This script re-aggregates the per-employee rows 11-18 (two rows per
# employee) down into a single summary row per employee (rows 11-14),
# then removes the now-unused trailing rows 15-18.
#
# New row 11 = Chrissy Cummings (old rows 17 + 18)
# New row 12 = Danielle Mai     (old rows 13 + 14)
# New row 13 = Jasmine Saiz     (old rows 15 + 16)
# New row 14 = Karen Trevizo    (old rows 11 + 12)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: Chrissy Cummings (sum of old rows 17 and 18) ---
$ws.Cells.Item(11, 1).Value = "Chrissy Cummings"
$ws.Cells.Item(11, 2).Value = 29
$ws.Cells.Item(11, 3).Value = 24
$ws.Cells.Item(11, 4).Value = 6
$ws.Cells.Item(11, 5).Value = 4.5
$ws.Cells.Item(11, 6).Value = 2094
$ws.Cells.Item(11, 7).Value = 2256.15
$ws.Cells.Item(11, 8).Value = 0
$ws.Cells.Item(11, 9).Value = 0
$ws.Cells.Item(11, 10).Value = 0
$ws.Cells.Item(11, 11).Value = 0
$ws.Cells.Item(11, 12).Value = 2094
$ws.Cells.Item(11, 13).Value = 2256.15
$ws.Cells.Item(11, 14).Value = 72.20689655172414
$ws.Cells.Item(11, 15).Value = 77.79827586206896

# --- Row 12: Danielle Mai (sum of old rows 13 and 14) ---
$ws.Cells.Item(12, 1).Value = "Danielle Mai"
$ws.Cells.Item(12, 2).Value = 9
$ws.Cells.Item(12, 3).Value = 7
$ws.Cells.Item(12, 4).Value = 1
$ws.Cells.Item(12, 5).Value = 0
$ws.Cells.Item(12, 6).Value = 571
$ws.Cells.Item(12, 7).Value = 614.88
$ws.Cells.Item(12, 8).Value = 0
$ws.Cells.Item(12, 9).Value = 0
$ws.Cells.Item(12, 10).Value = 30
$ws.Cells.Item(12, 11).Value = 32.31
$ws.Cells.Item(12, 12).Value = 601
$ws.Cells.Item(12, 13).Value = 647.1899999999999
$ws.Cells.Item(12, 14).Value = 66.77777777777777
$ws.Cells.Item(12, 15).Value = 71.91

# --- Row 13: Jasmine Saiz (sum of old rows 15 and 16) ---
$ws.Cells.Item(13, 1).Value = "Jasmine Saiz"
$ws.Cells.Item(13, 2).Value = 56
$ws.Cells.Item(13, 3).Value = 27
$ws.Cells.Item(13, 4).Value = 10
$ws.Cells.Item(13, 5).Value = 4.550000000000001
$ws.Cells.Item(13, 6).Value = 4107
$ws.Cells.Item(13, 7).Value = 4425.11
$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(13, 9).Value = 0
$ws.Cells.Item(13, 10).Value = 0
$ws.Cells.Item(13, 11).Value = 0
$ws.Cells.Item(13, 12).Value = 4107
$ws.Cells.Item(13, 13).Value = 4425.11
$ws.Cells.Item(13, 14).Value = 73.33928571428571
$ws.Cells.Item(13, 15).Value = 79.01982142857142

# --- Row 14: Karen Trevizo (sum of old rows 11 and 12) ---
$ws.Cells.Item(14, 1).Value = "Karen Trevizo"
$ws.Cells.Item(14, 2).Value = 27
$ws.Cells.Item(14, 3).Value = 21
$ws.Cells.Item(14, 4).Value = 1
$ws.Cells.Item(14, 5).Value = 4.699999999999999
$ws.Cells.Item(14, 6).Value = 1956
$ws.Cells.Item(14, 7).Value = 2106.72
$ws.Cells.Item(14, 8).Value = 0
$ws.Cells.Item(14, 9).Value = 0
$ws.Cells.Item(14, 10).Value = 14
$ws.Cells.Item(14, 11).Value = 15.08
$ws.Cells.Item(14, 12).Value = 1970
$ws.Cells.Item(14, 13).Value = 2121.8
$ws.Cells.Item(14, 14).Value = 72.96296296296296
$ws.Cells.Item(14, 15).Value = 78.5851851851852

# Remove old trailing rows 15-18, which are no longer needed now that
# their data has been folded into rows 11-14 above.
$ws.Range("A15:A18").EntireRow.Delete()
